$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Taher Haidari): clear arrival/departure date-time info
$ws.Range("C4:G4").ClearContents()

# Row 5 (Yousef Haidari): update arrival_time, departure_time, time_difference
$ws.Range("D5").Value = "20:28:02"
$ws.Range("F5").Value = "20:28:08"
$ws.Range("G5").Value = "0:00:06"
